# Auto-generated edit script applying numeric cell updates
# across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as described by the source diff (scheduled-runner profit recalculation).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 367.35715
$ws.Range("J19").Value = 426.4
$ws.Range("L19").Value = 426.4
$ws.Range("N19").Value = -776.4

# Row 100
$ws.Range("H100").Value = 18520208
$ws.Range("I100").Value = 23811182
$ws.Range("J100").Value = 1800
$ws.Range("K100").Value = 23811182
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -23810641
$ws.Range("N100").Value = -2882

# Row 103
$ws.Range("H103").Value = 271450.38
$ws.Range("I103").Value = 381.66666
$ws.Range("J103").Value = 361806.62
$ws.Range("K103").Value = 1144.99998
$ws.Range("L103").Value = 1085419.86
$ws.Range("M103").Value = -558.9999800000001
$ws.Range("N103").Value = -1086591.86

# Row 112
$ws.Range("H112").Value = 1429551.4
$ws.Range("I112").Value = 399.66666
$ws.Range("K112").Value = 1198.99998
$ws.Range("M112").Value = -90.99998000000005

# Row 129
$ws.Range("H129").Value = 1050.806
$ws.Range("I129").Value = 797.5
$ws.Range("J129").Value = 1066.8889
$ws.Range("K129").Value = 2392.5
$ws.Range("L129").Value = 3200.6667
$ws.Range("M129").Value = 2607.5
$ws.Range("N129").Value = -13200.6667

# Row 133
$ws.Range("H133").Value = 42890
$ws.Range("J133").Value = 42890
$ws.Range("L133").Value = 42890
$ws.Range("N133").Value = -53010

# Row 137
$ws.Range("H137").Value = 1746.7812
$ws.Range("I137").Value = 1587.88
$ws.Range("J137").Value = 2314.2856
$ws.Range("K137").Value = 4763.64
$ws.Range("L137").Value = 6942.8568
$ws.Range("M137").Value = -2213.64
$ws.Range("N137").Value = -12042.8568

# Row 138
$ws.Range("H138").Value = 1925.0492
$ws.Range("I138").Value = 1131
$ws.Range("J138").Value = 2440.1082
$ws.Range("K138").Value = 3393
$ws.Range("L138").Value = 7320.3246
$ws.Range("M138").Value = 1747
$ws.Range("N138").Value = -17600.3246

$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 86339.664
$ws.Range("I23").Value = 79506
$ws.Range("K23").Value = 79506
$ws.Range("M23").Value = -79247

# Row 45
$ws.Range("H45").Value = 12372.111
$ws.Range("I45").Value = 12372.111
$ws.Range("K45").Value = 12372.111
$ws.Range("M45").Value = -11995.111

# Row 74
$ws.Range("H74").Value = 1393.7059
$ws.Range("I74").Value = 1271.3636
$ws.Range("J74").Value = 1618
$ws.Range("K74").Value = 1271.3636
$ws.Range("L74").Value = 1618
$ws.Range("M74").Value = -397.3635999999999
$ws.Range("N74").Value = -3366

# Row 77
$ws.Range("H77").Value = 1393.7059
$ws.Range("I77").Value = 1271.3636
$ws.Range("J77").Value = 1618
$ws.Range("K77").Value = 6356.817999999999
$ws.Range("L77").Value = 8090
$ws.Range("M77").Value = -1988.817999999999
$ws.Range("N77").Value = -16826

# Row 102
$ws.Range("H102").Value = 2850832.2
$ws.Range("I102").Value = 2850832.2
$ws.Range("K102").Value = 2850832.2
$ws.Range("M102").Value = -2849210.2

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2291.0557
$ws.Range("I105").Value = 1939.9
$ws.Range("J105").Value = 2730
$ws.Range("K105").Value = 1939.9
$ws.Range("L105").Value = 2730
$ws.Range("M105").Value = -192.9000000000001
$ws.Range("N105").Value = -6224

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 333343330
$ws.Range("J23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15480

# Row 27
$ws.Range("H27").Value = 333343330
$ws.Range("J27").Value = 15000
$ws.Range("L27").Value = 15000
$ws.Range("N27").Value = -15384

# Row 31
$ws.Range("H31").Value = 3567.0842
$ws.Range("I31").Value = 1571.2094
$ws.Range("J31").Value = 5712.65
$ws.Range("K31").Value = 1571.2094
$ws.Range("L31").Value = 5712.65
$ws.Range("M31").Value = -1276.2094
$ws.Range("N31").Value = -6302.65

# Row 34
$ws.Range("H34").Value = 3567.0842
$ws.Range("I34").Value = 1571.2094
$ws.Range("J34").Value = 5712.65
$ws.Range("K34").Value = 1571.2094
$ws.Range("L34").Value = 5712.65
$ws.Range("M34").Value = -1369.2094
$ws.Range("N34").Value = -6116.65

$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 192.85715
$ws.Range("I18").Value = 208.33333
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 624.99999
$ws.Range("L18").Value = 300
$ws.Range("M18").Value = -455.99999
$ws.Range("N18").Value = -638

# Row 21
$ws.Range("H21").Value = 746.1539
$ws.Range("I21").Value = 300
$ws.Range("J21").Value = 880
$ws.Range("K21").Value = 900
$ws.Range("L21").Value = 2640
$ws.Range("M21").Value = -727
$ws.Range("N21").Value = -2986

# Row 118
$ws.Range("H118").Value = 2101.2856
$ws.Range("I118").Value = 903
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 2709
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = -1466
$ws.Range("N118").Value = -11486

$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 10751975
$ws.Range("I24").Value = 14333333
$ws.Range("J24").Value = 7900
$ws.Range("K24").Value = 14333333
$ws.Range("L24").Value = 7900
$ws.Range("M24").Value = -14333160
$ws.Range("N24").Value = -8246

# Row 70
$ws.Range("H70").Value = 5362.485
$ws.Range("I70").Value = 5315
$ws.Range("J70").Value = 5510.875
$ws.Range("K70").Value = 5315
$ws.Range("L70").Value = 5510.875
$ws.Range("M70").Value = -5045
$ws.Range("N70").Value = -6050.875

# Row 73
$ws.Range("H73").Value = 5362.485
$ws.Range("I73").Value = 5315
$ws.Range("J73").Value = 5510.875
$ws.Range("K73").Value = 5315
$ws.Range("L73").Value = 5510.875
$ws.Range("M73").Value = -4379
$ws.Range("N73").Value = -7382.875

# Row 97
$ws.Range("H97").Value = 2424.8333
$ws.Range("I97").Value = 2424.8333
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2424.8333
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1928.8333
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2550
$ws.Range("I93").Value = 3900.3333
$ws.Range("J93").Value = 1199.6666
$ws.Range("K93").Value = 3900.3333
$ws.Range("L93").Value = 1199.6666
$ws.Range("M93").Value = -2652.3333
$ws.Range("N93").Value = -3695.6666

# Row 104
$ws.Range("H104").Value = 11122.714
$ws.Range("J104").Value = 11122.714
$ws.Range("L104").Value = 11122.714
$ws.Range("N104").Value = -18110.714

$ws = $wb.Worksheets.Item("WVR")
# Row 123
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

